$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 84, shifting existing rows 84:139 down to 85:140
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new weekly data entry
$ws.Cells.Item(84, 1).Value = 10
$ws.Cells.Item(84, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(84, 3).Value = 'La Araucanía'
$ws.Cells.Item(84, 4).Value = 44806
$ws.Cells.Item(84, 5).Value = 9
$ws.Cells.Item(84, 6).Value = 'Fruta'
$ws.Cells.Item(84, 7).Value = 100107
$ws.Cells.Item(84, 8).Value = 'Otros'
$ws.Cells.Item(84, 9).Value = 100107002
$ws.Cells.Item(84, 10).Value = 'Chirimoya'
$ws.Cells.Item(84, 11).Value = 'Cultivar IV Región'
$ws.Cells.Item(84, 12).Value = 'Primera'
$ws.Cells.Item(84, 13).Value = 50
$ws.Cells.Item(84, 14).Value = 3000
$ws.Cells.Item(84, 15).Value = 3000
$ws.Cells.Item(84, 16).Value = 3000
$ws.Cells.Item(84, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(84, 18).Value = 'Provincia del Elquí'
$ws.Cells.Item(84, 19).Value = 3000
$ws.Cells.Item(84, 20).Value = 1
